$d = $word.ActiveDocument

# The "Requirement Name / Priority" overview table (second table in the
# document) ended up with seven stray trailing rows duplicating/extending
# the S1-S7 summary (an "S7. Programming language" row plus six detail
# rows reusing the S1-S6 labels with longer descriptions). Remove them so
# the table goes back to just its header + the original S1-S6 rows.
$t = $d.Tables.Item(2)

$keep = 7
for ($i = $t.Rows.Count; $i -gt $keep; $i--) {
    $t.Rows.Item($i).Delete()
}

Write-Output ("Table2 rows now: " + $t.Rows.Count)
